$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Delete old rows 59:65 (sofa_max/respiratory_max/coagulation_max/liver_max/cardiovascular_max/cns_max/renal_max_0_24h rows removed) ---
# This shifts old rows 66-91 up by 7 to become new rows 59-84, and the sheet dimension becomes A1:E84
$ws.Rows("59:65").Delete() | Out-Null

# --- Step 2: Update values in rows 4-58 (row numbers unchanged by the deletion above) ---
$ws.Range("D4").Value = 1467
$ws.Range("E4").Value = 6130
$ws.Range("D5").Value = '612 (41.7)'
$ws.Range("E5").Value = '2291 (37.4)'
$ws.Range("D6").Value = '350 (23.9)'
$ws.Range("E6").Value = '1411 (23.0)'
$ws.Range("D7").Value = '298 (20.3)'
$ws.Range("E7").Value = '1385 (22.6)'
$ws.Range("D8").Value = '207 (14.1)'
$ws.Range("E8").Value = '1043 (17.0)'
$ws.Range("D9").Value = '687 (46.8)'
$ws.Range("E9").Value = '2591 (42.3)'
$ws.Range("D10").Value = '557 (38.0)'
$ws.Range("E10").Value = '3014 (49.2)'
$ws.Range("D11").Value = '204 (13.9)'
$ws.Range("E11").Value = '348 (5.7)'
$ws.Range("D12").Value = '706 (48.1)'
$ws.Range("E12").Value = '2768 (45.2)'
$ws.Range("D13").Value = '461 (31.4)'
$ws.Range("E13").Value = '308 (5.0)'
$ws.Range("D14").Value = '101 (6.9)'
$ws.Range("E14").Value = '718 (11.7)'
$ws.Range("D15").Value = '511 (34.8)'
$ws.Range("E15").Value = '2592 (42.3)'
$ws.Range("D16").Value = '303 (20.7)'
$ws.Range("E16").Value = '1222 (19.9)'
$ws.Range("D17").Value = '433 (29.5)'
$ws.Range("E17").Value = '1776 (29.0)'
$ws.Range("D18").Value = '957 (65.2)'
$ws.Range("E18").Value = '4204 (68.6)'
$ws.Range("D19").Value = '216 (14.7)'
$ws.Range("E19").Value = '604 (9.9)'
$ws.Range("D20").Value = '745 (50.8)'
$ws.Range("E20").Value = '3520 (57.4)'
$ws.Range("D21").Value = '1036 (70.6)'
$ws.Range("E21").Value = '4508 (73.5)'
$ws.Range("D22").Value = '280 (19.1)'
$ws.Range("E22").Value = '829 (13.5)'
$ws.Range("D23").Value = '830 (56.6)'
$ws.Range("E23").Value = '3860 (63.0)'
$ws.Range("D24").Value = '995 (67.8)'
$ws.Range("E24").Value = '3993 (65.1)'
$ws.Range("D25").Value = '178 (12.1)'
$ws.Range("E25").Value = '879 (14.3)'
$ws.Range("D26").Value = '1449 (98.8)'
$ws.Range("E26").Value = '6078 (99.2)'
$ws.Range("D27").Value = '1044 (71.2)'
$ws.Range("E27").Value = '3970 (64.8)'
$ws.Range("D28").Value = '585 (39.9)'
$ws.Range("E28").Value = '2348 (38.3)'
$ws.Range("D29").Value = '362 (24.7)'
$ws.Range("E29").Value = '1567 (25.6)'
$ws.Range("D30").Value = '20 (1.4)'
$ws.Range("E30").Value = '103 (1.7)'
$ws.Range("D31").Value = '436 (29.7)'
$ws.Range("E31").Value = '2268 (37.0)'
$ws.Range("D32").Value = '1 (0.1)'
$ws.Range("D33").Value = '9 (0.6)'
$ws.Range("E33").Value = '42 (0.7)'
$ws.Range("D34").Value = '114 (7.8)'
$ws.Range("E34").Value = '296 (4.8)'
$ws.Range("D35").Value = '41 (2.8)'
$ws.Range("E35").Value = '114 (1.9)'
$ws.Range("D36").Value = '169 (11.5)'
$ws.Range("E36").Value = '307 (5.0)'
$ws.Range("D37").Value = '1133 (77.2)'
$ws.Range("E37").Value = '5371 (87.6)'
$ws.Range("D38").Value = '48 (3.3)'
$ws.Range("E38").Value = '164 (2.7)'
$ws.Range("D39").Value = '597 (40.7)'
$ws.Range("E39").Value = '1709 (27.9)'
$ws.Range("D40").Value = '822 (56.0)'
$ws.Range("E40").Value = '4257 (69.4)'
$ws.Range("D41").Value = '70 (4.8)'
$ws.Range("E41").Value = '289 (4.7)'
$ws.Range("D42").Value = '83 (5.7)'
$ws.Range("E42").Value = '339 (5.5)'
$ws.Range("D43").Value = '11 (0.7)'
$ws.Range("E43").Value = '38 (0.6)'
$ws.Range("D44").Value = ""
$ws.Range("E44").Value = '11 (0.2)'
$ws.Range("D45").Value = '2 (0.1)'
$ws.Range("D46").Value = '64 [51,74]'
$ws.Range("E46").Value = '67 [57,77]'
$ws.Range("C47").Value = 6072
$ws.Range("D47").Value = '7.63 [5.44,11.84]'
$ws.Range("E47").Value = '7.88 [5.54,11.83]'
$ws.Range("C48").Value = 1525
$ws.Range("D48").Value = '7.21 [5.08,11.83]'
$ws.Range("E48").Value = '7.04 [5.04,11.01]'
$ws.Range("C49").Value = 6072
$ws.Range("D49").Value = '13.00 [8.00,22.00]'
$ws.Range("C50").Value = 1525
$ws.Range("D50").Value = '16.00 [11.00,25.00]'
$ws.Range("E50").Value = '15.00 [10.00,23.00]'
$ws.Range("C53").Value = 2265
$ws.Range("D53").Value = '2 [1,4]'
$ws.Range("C54").Value = 14
$ws.Range("C55").Value = 2789
$ws.Range("C56").Value = 17
$ws.Range("C57").Value = 21

# --- Step 3: Update values in new rows 59-84 (previously old rows 66-91, shifted up by the deletion) ---
$ws.Range("C59").Value = 2844
$ws.Range("D59").Value = '900 [279,1569]'
$ws.Range("E59").Value = '848 [275,1660]'
$ws.Range("C60").Value = 70
$ws.Range("D60").Value = '4405 [2011,8575]'
$ws.Range("E60").Value = '4132 [1773,8008]'
$ws.Range("C61").Value = 70
$ws.Range("D61").Value = '532.7 [276.3,912.8]'
$ws.Range("E61").Value = '505.5 [255.0,870.4]'
$ws.Range("C62").Value = 4081
$ws.Range("D62").Value = '50 [40,61]'
$ws.Range("E62").Value = '50 [42,65]'
$ws.Range("C63").Value = 2053
$ws.Range("D63").Value = '60.0 [25.0,115.0]'
$ws.Range("E63").Value = '51.0 [21.0,104.0]'
$ws.Range("C64").Value = 2053
$ws.Range("D64").Value = '0.30 [0.14,0.50]'
$ws.Range("E64").Value = '0.26 [0.12,0.46]'
$ws.Range("C65").Value = 2053
$ws.Range("D65").Value = '2.5 [1.0,12.0]'
$ws.Range("E65").Value = '3.0 [1.0,13.0]'
$ws.Range("C66").Value = 6488
$ws.Range("D66").Value = '22.0 [4.7,62.5]'
$ws.Range("E66").Value = '32.0 [7.0,76.0]'
$ws.Range("C67").Value = 2858
$ws.Range("D67").Value = '4.0 [1.0,22.0]'
$ws.Range("E67").Value = '4.0 [1.0,19.0]'
$ws.Range("C68").Value = 2858
$ws.Range("D68").Value = '45.0 [14.0,94.0]'
$ws.Range("E68").Value = '44.0 [17.0,92.0]'
$ws.Range("C69").Value = 2858
$ws.Range("D69").Value = '0.21 [0.07,0.46]'
$ws.Range("E69").Value = '0.24 [0.08,0.46]'
$ws.Range("C70").Value = 22
$ws.Range("D70").Value = '19.7 [17.1,22.9]'
$ws.Range("E70").Value = '19.4 [17.1,22.4]'
$ws.Range("C71").Value = 17
$ws.Range("D71").Value = '77.9 [71.5,86.4]'
$ws.Range("E71").Value = '75.1 [69.5,82.3]'
$ws.Range("C72").Value = 430
$ws.Range("D72").Value = '36.9 [36.6,37.4]'
$ws.Range("E72").Value = '36.9 [36.6,37.3]'
$ws.Range("C73").Value = 18
$ws.Range("D73").Value = '97.9 [96.2,99.2]'
$ws.Range("E73").Value = '97.3 [95.7,98.6]'
$ws.Range("C74").Value = 17
$ws.Range("D74").Value = '88.4 [76.6,101.4]'
$ws.Range("E74").Value = '87.1 [76.1,99.3]'
$ws.Range("C75").Value = 1917
$ws.Range("D75").Value = '84.5 [66.0,119.0]'
$ws.Range("E75").Value = '85.0 [69.0,113.0]'
$ws.Range("C76").Value = 1917
$ws.Range("D76").Value = '44.0 [37.0,53.0]'
$ws.Range("E76").Value = '46.0 [39.0,54.0]'
$ws.Range("C77").Value = 1083
$ws.Range("D77").Value = '7.3 [7.2,7.4]'
$ws.Range("E77").Value = '7.3 [7.2,7.4]'
$ws.Range("C78").Value = 28
$ws.Range("D78").Value = '160.0 [124.0,226.0]'
$ws.Range("E78").Value = '152.0 [122.0,200.0]'
$ws.Range("C79").Value = 11
$ws.Range("D79").Value = '137.0 [134.0,140.0]'
$ws.Range("E79").Value = '137.0 [134.0,140.0]'
$ws.Range("C80").Value = 15
$ws.Range("D80").Value = '4.5 [4.1,5.2]'
$ws.Range("E80").Value = '4.5 [4.1,5.0]'
$ws.Range("C81").Value = 7411
$ws.Range("D81").Value = '17.8 [13.0,28.8]'
$ws.Range("E81").Value = '23.4 [13.5,36.8]'
$ws.Range("C82").Value = 1022
$ws.Range("D82").Value = '9.9 [8.3,11.6]'
$ws.Range("E82").Value = '10.0 [8.5,11.6]'
$ws.Range("C83").Value = 4934
$ws.Range("D83").Value = '228.0 [150.5,365.0]'
$ws.Range("E83").Value = '235.0 [160.0,359.5]'
$ws.Range("C84").Value = 441
$ws.Range("D84").Value = '1.4 [1.2,1.8]'
$ws.Range("E84").Value = '1.4 [1.2,1.8]'
